# Update "想去人数" (column F) counts across the workbook to the freshly
# regenerated values (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# sheet name -> @{ row = newValue }
$updates = @{
    "展览" = @{
        4  = 504
        5  = 1187
        7  = 202
        8  = 72
        9  = 781
        10 = 430
        11 = 59
        12 = 271
        14 = 207
        17 = 6325
        19 = 61
        21 = 7295
        24 = 3315
        25 = 440
        26 = 819
        27 = 4483
        29 = 167
        30 = 166
        31 = 1334
        32 = 127
        35 = 1047
        36 = 1401
        37 = 2092
    }
    "本地生活" = @{
        3 = 1179
    }
    "全部类型" = @{
        4  = 1179
        7  = 504
        8  = 1187
        10 = 202
        11 = 72
        12 = 781
        13 = 430
        14 = 59
        15 = 271
        18 = 207
        21 = 6325
        23 = 61
        25 = 7295
        28 = 3315
        29 = 440
        30 = 819
        31 = 4483
        34 = 167
        35 = 166
        36 = 1334
        37 = 127
        40 = 1047
        41 = 1401
        43 = 2092
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
